# Generate Report for Handback
#
# Updates the handoff/handback timestamps for the "f9878666-..." file
# (row 3 of each per-language sheet) and the corresponding "Latest HO
# Xliff Generate Date" on the Overview sheet, reflecting a freshly
# generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" for f9878666-90b1-48cc-b4fa-4a6a87ee1180.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-25 22:47:31"

# --- zh-cn sheet ---
# Row 3 corresponds to f9878666-90b1-48cc-b4fa-4a6a87ee1180
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-25 22:47:27"   # Correspond Handoff Datetime
$wsZhCn.Range("K3").Value = "2016-08-25 22:47:48"   # Correspond Handback DateTime

# --- de-de sheet ---
# Row 3 corresponds to f9878666-90b1-48cc-b4fa-4a6a87ee1180
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-25 22:47:31"   # Correspond Handoff Datetime
$wsDeDe.Range("K3").Value = "2016-08-25 22:47:54"   # Correspond Handback DateTime
